# Fruta / hortaliza, semanal
#
# Inserts a new weekly price-report row for "Feria Lagunitas de Puerto Montt"
# (Santina cherries, week of 2022-01-07) ahead of the existing row 56,
# pushing the subsequent rows (old 56-59) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 56..59 down to 57..60, leaving row 56 free for the
# new record.
$ws.Rows(56).Insert()

# Populate the newly inserted row 56 with the new weekly record.
$ws.Range("A56").Value = 4
$ws.Range("B56").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C56").Value = "Los Lagos"
$ws.Range("D56").Value = 44568
$ws.Range("E56").Value = 10
$ws.Range("F56").Value = "Fruta"
$ws.Range("G56").Value = 100103
$ws.Range("H56").Value = "Frutos de hueso (carozo)"
$ws.Range("I56").Value = 100103001
$ws.Range("J56").Value = "Cereza"
$ws.Range("K56").Value = "Santina"
$ws.Range("L56").Value = "Primera"
$ws.Range("M56").Value = 800
$ws.Range("N56").Value = 7500
$ws.Range("O56").Value = 8000
$ws.Range("P56").Value = 7750
$ws.Range("Q56").Value = '$/bandeja 10 kilos'
$ws.Range("R56").Value = "Provincia de Curicó"
$ws.Range("S56").Value = 775
$ws.Range("T56").Value = 10
